# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.851.83"
$ws.Range("E2").Value = "  -3.11%  "

$ws.Range("D3").Value = "1.620.01"
$ws.Range("E3").Value = "  -3.11%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'308.39"
$ws.Range("E5").Value = "  -1.74%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "'0.3934"
$ws.Range("E7").Value = "  -0.49%  "

$ws.Range("D8").Value = "'0.3840"
$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("D9").Value = "'1.000"
$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.361"
$ws.Range("E10").Value = "  -2.50%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'49.41"
$ws.Range("E11").Value = "  -2.66%  "

$ws.Range("D12").Value = "'0.08458"
$ws.Range("E12").Value = "  -2.10%  "

$ws.Range("E13").Value = "  -6.35%  "

$ws.Range("D14").Value = "'7.056"
$ws.Range("E14").Value = "  -3.43%  "

$ws.Range("D15").Value = "'7.587"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "'0.00001280"
$ws.Range("E16").Value = "  -2.98%  "

$ws.Range("D17").Value = "1.624.44"
$ws.Range("E17").Value = "  -4.02%  "

$ws.Range("D18").Value = "'93.81"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").Value = "'0.06935"
$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("D20").Value = "'20.06"
$ws.Range("E20").Value = "  -4.91%  "

$ws.Range("D21").Value = "'6.822"
$ws.Range("E21").Value = "  -3.58%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'13.43"
$ws.Range("E23").Value = "  -3.48%  "

$ws.Range("D24").Value = "23.857.68"
$ws.Range("E24").Value = "  -3.05%  "

$ws.Range("D25").Value = "'2.457"
$ws.Range("E25").Value = "  +4.75%  "

$ws.Range("D26").Value = "'2.853"
$ws.Range("E26").Value = "  +3.26%  "

$ws.Range("D27").Value = "'22.25"
$ws.Range("E27").Value = "  -3.27%  "

$ws.Range("D28").Value = "'157.03"
$ws.Range("E28").Value = "  -1.95%  "

$ws.Range("D29").Value = "'139.90"
$ws.Range("E29").Value = "  -3.92%  "

$ws.Range("D30").Value = "'5.273"
$ws.Range("E30").Value = "  -9.73%  "

$ws.Range("D31").Value = "'7.845"
$ws.Range("E31").Value = "  -5.39%  "

$ws.Range("D32").Value = "'2.489"
$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("D33").Value = "1.793.70"
$ws.Range("E33").Value = "  -3.64%  "

$ws.Range("D34").Value = "'0.08102"
$ws.Range("E34").Value = "  -1.96%  "

$ws.Range("D35").Value = "'0.9726"
$ws.Range("E35").Value = "  -2.25%  "

$ws.Range("D36").Value = "'0.02894"
$ws.Range("E36").Value = "  -5.95%  "

$ws.Range("D37").Value = "'6.584"
$ws.Range("E37").Value = "  -4.69%  "

$ws.Range("D38").Value = "'0.2667"
$ws.Range("E38").Value = "  -4.90%  "

$ws.Range("D39").Value = "'0.09144"
$ws.Range("E39").Value = "  -4.99%  "

$ws.Range("D40").Value = "'10.36"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").Value = "'13.61"
$ws.Range("E41").Value = "  +0.42%  "

$ws.Range("D42").Value = "'1.432"
$ws.Range("E42").Value = "  -5.70%  "

$ws.Range("D43").Value = "'0.7504"
$ws.Range("E43").Value = "  -4.91%  "

$ws.Range("D44").Value = "'16.14"
$ws.Range("E44").Value = "  -2.41%  "

$ws.Range("D45").Value = "'0.6914"
$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("D46").Value = "'2.471"
$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").Value = "'4.076"
$ws.Range("E47").Value = "  -2.31%  "

$ws.Range("D48").Value = "'0.9998"
$ws.Range("E48").Value = "  -0.19%  "

$ws.Range("D49").Value = "'0.08237"
$ws.Range("E49").Value = "  -4.53%  "

$ws.Range("D50").Value = "'134.38"
$ws.Range("E50").Value = "  -2.84%  "

$ws.Range("D51").Value = "'1.205"
$ws.Range("E51").Value = "  -9.05%  "
